# Auto-generated edit script: update crypto Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.089.10"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.833.12"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'243.34"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "'0.6274"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "'0.07499"
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("D9").Value = "'0.2919"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").Value = "'23.18"
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("D11").Value = "'0.07685"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "1.832.57"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "'0.6668"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "'82.70"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("E16").Value = "  -8.30%  "
$ws.Range("D17").Value = "'5.975"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").Value = "29.104.24"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").Value = "2.081.89"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'12.57"
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("D21").Value = "'223.21"
$ws.Range("E21").Value = "  -1.47%  "
$ws.Range("D22").Value = "'1.003"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").Value = "'7.095"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").Value = "'159.75"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").Value = "'0.1393"
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("D27").Value = "'8.491"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'17.91"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "'1.499"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").Value = "'0.05685"
$ws.Range("E30").Value = "  +8.73%  "
$ws.Range("D31").Value = "'4.148"
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("D32").Value = "'4.088"
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "'0.7425"
$ws.Range("E35").Value = "  +0.98%  "
$ws.Range("D36").Value = "'1.140"
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").Value = "'2.762"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").Value = "1.222.12"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").Value = "'0.01779"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").Value = "'6.536"
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").Value = "'101.92"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").Value = "1.982.57"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "'65.80"
$ws.Range("E47").Value = "  +2.68%  "
$ws.Range("D48").Value = "'0.5088"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("D50").Value = "'0.07431"
$ws.Range("E50").Value = "  +6.50%  "
$ws.Range("D51").Value = "'8.978"
$ws.Range("E51").Value = "  +0.53%  "
